# Remove the trailing copyright block that was dropped from the page:
#   - an empty "Normal" paragraph
#   - an empty "Normal" paragraph with PageBreakBefore
#   - the "© 2020 . Contact: ..." paragraph
# while leaving the final two (still-empty) paragraphs and the section
# properties untouched.

$d = $word.ActiveDocument

# Locate the paragraph that carries the copyright notice.
$copyright = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Contact: luizeleno@usp.br*") {
        $copyright = $p
        break
    }
}

if ($copyright -eq $null) {
    throw "Could not locate the copyright paragraph"
}

# The paragraph immediately before it is the empty page-break paragraph,
# and the one before that is the plain empty paragraph - both are being
# removed together with the copyright paragraph itself.
$pageBreakPara = $copyright.Previous(1)
$blankPara = $pageBreakPara.Previous(1)

$deleteRange = $d.Range($blankPara.Range.Start, $copyright.Range.End)
$deleteRange.Delete()
